$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "codice _1_livello"
$ws.Cells.Item(1,2).Value = "label_ITA_1_livello"
$ws.Cells.Item(1,3).Value = "label_ENG_1_livello"
$ws.Cells.Item(1,4).Value = "label_DEU_1_livello"

# Data rows: Code, Italian label, English label, German label
$data = @(
    @("REQ",        "Istanza/richiesta",              "Request/application",           "Gesuch / Anfrage"),
    @("OTHDOC",      "Altra documentazione ",           "Other documents",               "Sonstige Dokumentation"),
    @("PAYMENTDEC",  "Attestazione di pagamento",        "Payment declaration",           "Zahlungsbestätigung"),
    @("AUTHACT",     "Atto autorizzativo",               "Authorization Act",             "Bewilligungsurkunde"),
    @("IDDEC",       "Attestazione di identità ",        "Identity Declaration",          "Identitätsnachweis"),
    @("ADMINDOC",    "Documentazione amministrativa",    "Administrative documentation",  "Verwaltungsdokumentation"),
    @("CERT",        "Certificazione",                   "Certification",                 "Bescheinigung"),
    @("CODE",        "Codice",                           "Code",                          "Kode")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

# Style: header row bold (already existing style carries over since B1 was bold); apply bold to new C1/D1 headers
$ws.Range("A1:D1").Font.Bold = $true

# Column D (German translations) uses a distinct font (no explicit color) - apply font name/size to match
$ws.Range("D2:D9").Font.Name = "Calibri"
$ws.Range("D2:D9").Font.Size = 12

# Column widths (columns A and B already have the correct width from the
# source workbook; only the two new columns need sizing - these were
# best-fit/auto-fitted to their (translated) content in the real edit)
$ws.Columns.Item(3).ColumnWidth = 25.330729166666668
$ws.Columns.Item(4).ColumnWidth = 23.166666666666668

# Selection matches the diff (activeCell D2)
$ws.Range("D2").Select()
